# Update "想去人数" (interest count, column F) and occasionally
# "最低票价" (minimum price, column G) figures on the "展览" and
# "全部类型" worksheets, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 5519
$ws1.Range("F3").Value  = 110
$ws1.Range("F5").Value  = 354
$ws1.Range("F6").Value  = 47
$ws1.Range("F8").Value  = 20
$ws1.Range("F9").Value  = 28
$ws1.Range("F10").Value = 57
$ws1.Range("F13").Value = 336
$ws1.Range("F14").Value = 432
$ws1.Range("F15").Value = 3031
$ws1.Range("F16").Value = 2
$ws1.Range("G16").Value = 45
$ws1.Range("F18").Value = 1646
$ws1.Range("F19").Value = 15

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 5519
$ws4.Range("F3").Value  = 110
$ws4.Range("F5").Value  = 354
$ws4.Range("F6").Value  = 47
$ws4.Range("F9").Value  = 20
$ws4.Range("F10").Value = 28
$ws4.Range("F11").Value = 57
$ws4.Range("F14").Value = 336
$ws4.Range("F15").Value = 432
$ws4.Range("F16").Value = 3031
$ws4.Range("F17").Value = 2
$ws4.Range("G17").Value = 45
$ws4.Range("F19").Value = 1646
$ws4.Range("F20").Value = 15
